# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells, 1 column past the existing "Unnamed: 28" (AC) column:
# AD1 = Wins, AE1 = Losses, AF1 = Ties.
# Copy the header style (bold, bordered, centered) from A1 so the new
# headers match the rest of row 1.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record for every player row (2-64) with the team's
# overall record for the season: 91 wins, 71 losses, 0 ties.
for ($r = 2; $r -le 64; $r++) {
    $ws.Cells.Item($r, 30).Value = 91
    $ws.Cells.Item($r, 31).Value = 71
    $ws.Cells.Item($r, 32).Value = 0
}
